# Update the LR-pairs data table with refreshed TPM-based NATMI results.
# The new run adds the "Inflammatory-Mac" target cluster (column D) for every
# sending cluster, growing the table from 15 to 20 data rows (rows 2-21).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Nlgn2"
$ws.Cells.Item(2,3).Value = "Nrxn2"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 4.682284333333333
$ws.Cells.Item(2,8).Value = 14.046853
$ws.Cells.Item(2,9).Value = 0.1747450949072225
$ws.Cells.Item(2,10).Value = 0.1747450949072225
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 0.173461
$ws.Cells.Item(2,14).Value = 0.520383
$ws.Cells.Item(2,15).Value = 0.2181774959583017
$ws.Cells.Item(2,16).Value = 0.2181774959583018
$ws.Cells.Item(2,17).Value = 0.8121937227443332
$ws.Cells.Item(2,18).Value = 7.309743504699
$ws.Cells.Item(2,19).Value = 0.03812544723785359
$ws.Cells.Item(2,20).Value = 0.0381254472378536

# Row 3: ECs -> FAPs
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Nlgn2"
$ws.Cells.Item(3,3).Value = "Nrxn2"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 4.682284333333333
$ws.Cells.Item(3,8).Value = 14.046853
$ws.Cells.Item(3,9).Value = 0.1747450949072225
$ws.Cells.Item(3,10).Value = 0.1747450949072225
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 0.4903776666666667
$ws.Cells.Item(3,14).Value = 1.471133
$ws.Cells.Item(3,15).Value = 0.6167920822963554
$ws.Cells.Item(3,16).Value = 0.6167920822963555
$ws.Cells.Item(3,17).Value = 2.296087666049889
$ws.Cells.Item(3,18).Value = 20.664788994449
$ws.Cells.Item(3,19).Value = 0.1077813909589
$ws.Cells.Item(3,20).Value = 0.1077813909589

# Row 4: ECs -> Inflammatory-Mac
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Nlgn2"
$ws.Cells.Item(4,3).Value = "Nrxn2"
$ws.Cells.Item(4,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 4.682284333333333
$ws.Cells.Item(4,8).Value = 14.046853
$ws.Cells.Item(4,9).Value = 0.1747450949072225
$ws.Cells.Item(4,10).Value = 0.1747450949072225
$ws.Cells.Item(4,11).Value = 1
$ws.Cells.Item(4,12).Value = 0.3333333333333333
$ws.Cells.Item(4,13).Value = 0.04417666666666667
$ws.Cells.Item(4,14).Value = 0.13253
$ws.Cells.Item(4,15).Value = 0.0555649656874912
$ws.Cells.Item(4,16).Value = 0.0555649656874912
$ws.Cells.Item(4,17).Value = 0.2068477142322222
$ws.Cells.Item(4,18).Value = 1.86162942809
$ws.Cells.Item(4,19).Value = 0.009709705202577211
$ws.Cells.Item(4,20).Value = 0.009709705202577213

# Row 5: ECs -> MuSCs
$ws.Cells.Item(5,1).Value = "ECs"
$ws.Cells.Item(5,2).Value = "Nlgn2"
$ws.Cells.Item(5,3).Value = "Nrxn2"
$ws.Cells.Item(5,4).Value = "MuSCs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 4.682284333333333
$ws.Cells.Item(5,8).Value = 14.046853
$ws.Cells.Item(5,9).Value = 0.1747450949072225
$ws.Cells.Item(5,10).Value = 0.1747450949072225
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 0.08703
$ws.Cells.Item(5,14).Value = 0.26109
$ws.Cells.Item(5,15).Value = 0.1094654560578516
$ws.Cells.Item(5,16).Value = 0.1094654560578516
$ws.Cells.Item(5,17).Value = 0.4074992055299999
$ws.Cells.Item(5,18).Value = 3.667492849769999
$ws.Cells.Item(5,19).Value = 0.01912855150789168
$ws.Cells.Item(5,20).Value = 0.01912855150789168

# Row 6: FAPs -> ECs
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Nlgn2"
$ws.Cells.Item(6,3).Value = "Nrxn2"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 17.62337166666667
$ws.Cells.Item(6,8).Value = 52.870115
$ws.Cells.Item(6,9).Value = 0.6577126751045782
$ws.Cells.Item(6,10).Value = 0.6577126751045781
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 0.173461
$ws.Cells.Item(6,14).Value = 0.520383
$ws.Cells.Item(6,15).Value = 0.2181774959583017
$ws.Cells.Item(6,16).Value = 0.2181774959583018
$ws.Cells.Item(6,17).Value = 3.056967672671667
$ws.Cells.Item(6,18).Value = 27.512709054045
$ws.Cells.Item(6,19).Value = 0.1434981045143529
$ws.Cells.Item(6,20).Value = 0.143498104514353

# Row 7: FAPs -> FAPs
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Nlgn2"
$ws.Cells.Item(7,3).Value = "Nrxn2"
$ws.Cells.Item(7,4).Value = "FAPs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 17.62337166666667
$ws.Cells.Item(7,8).Value = 52.870115
$ws.Cells.Item(7,9).Value = 0.6577126751045782
$ws.Cells.Item(7,10).Value = 0.6577126751045781
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 0.4903776666666667
$ws.Cells.Item(7,14).Value = 1.471133
$ws.Cells.Item(7,15).Value = 0.6167920822963554
$ws.Cells.Item(7,16).Value = 0.6167920822963555
$ws.Cells.Item(7,17).Value = 8.642107876699445
$ws.Cells.Item(7,18).Value = 77.778970890295
$ws.Cells.Item(7,19).Value = 0.4056719704304591
$ws.Cells.Item(7,20).Value = 0.4056719704304591

# Row 8: FAPs -> Inflammatory-Mac
$ws.Cells.Item(8,1).Value = "FAPs"
$ws.Cells.Item(8,2).Value = "Nlgn2"
$ws.Cells.Item(8,3).Value = "Nrxn2"
$ws.Cells.Item(8,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 17.62337166666667
$ws.Cells.Item(8,8).Value = 52.870115
$ws.Cells.Item(8,9).Value = 0.6577126751045782
$ws.Cells.Item(8,10).Value = 0.6577126751045781
$ws.Cells.Item(8,11).Value = 1
$ws.Cells.Item(8,12).Value = 0.3333333333333333
$ws.Cells.Item(8,13).Value = 0.04417666666666667
$ws.Cells.Item(8,14).Value = 0.13253
$ws.Cells.Item(8,15).Value = 0.0555649656874912
$ws.Cells.Item(8,16).Value = 0.0555649656874912
$ws.Cells.Item(8,17).Value = 0.7785418156611112
$ws.Cells.Item(8,18).Value = 7.00687634095
$ws.Cells.Item(8,19).Value = 0.03654578222441393
$ws.Cells.Item(8,20).Value = 0.03654578222441393

# Row 9: FAPs -> MuSCs
$ws.Cells.Item(9,1).Value = "FAPs"
$ws.Cells.Item(9,2).Value = "Nlgn2"
$ws.Cells.Item(9,3).Value = "Nrxn2"
$ws.Cells.Item(9,4).Value = "MuSCs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 17.62337166666667
$ws.Cells.Item(9,8).Value = 52.870115
$ws.Cells.Item(9,9).Value = 0.6577126751045782
$ws.Cells.Item(9,10).Value = 0.6577126751045781
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 0.08703
$ws.Cells.Item(9,14).Value = 0.26109
$ws.Cells.Item(9,15).Value = 0.1094654560578516
$ws.Cells.Item(9,16).Value = 0.1094654560578516
$ws.Cells.Item(9,17).Value = 1.53376203615
$ws.Cells.Item(9,18).Value = 13.80385832535
$ws.Cells.Item(9,19).Value = 0.07199681793535224
$ws.Cells.Item(9,20).Value = 0.07199681793535224

# Row 10: Inflammatory-Mac -> ECs
$ws.Cells.Item(10,1).Value = "Inflammatory-Mac"
$ws.Cells.Item(10,2).Value = "Nlgn2"
$ws.Cells.Item(10,3).Value = "Nrxn2"
$ws.Cells.Item(10,4).Value = "ECs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 0.9374000000000001
$ws.Cells.Item(10,8).Value = 2.8122
$ws.Cells.Item(10,9).Value = 0.03498421716936109
$ws.Cells.Item(10,10).Value = 0.03498421716936108
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 0.173461
$ws.Cells.Item(10,14).Value = 0.520383
$ws.Cells.Item(10,15).Value = 0.2181774959583017
$ws.Cells.Item(10,16).Value = 0.2181774959583018
$ws.Cells.Item(10,17).Value = 0.1626023414
$ws.Cells.Item(10,18).Value = 1.4634210726
$ws.Cells.Item(10,19).Value = 0.007632768900072629
$ws.Cells.Item(10,20).Value = 0.007632768900072629

# Row 11: Inflammatory-Mac -> FAPs
$ws.Cells.Item(11,1).Value = "Inflammatory-Mac"
$ws.Cells.Item(11,2).Value = "Nlgn2"
$ws.Cells.Item(11,3).Value = "Nrxn2"
$ws.Cells.Item(11,4).Value = "FAPs"
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 0.9374000000000001
$ws.Cells.Item(11,8).Value = 2.8122
$ws.Cells.Item(11,9).Value = 0.03498421716936109
$ws.Cells.Item(11,10).Value = 0.03498421716936108
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 0.4903776666666667
$ws.Cells.Item(11,14).Value = 1.471133
$ws.Cells.Item(11,15).Value = 0.6167920822963554
$ws.Cells.Item(11,16).Value = 0.6167920822963555
$ws.Cells.Item(11,17).Value = 0.4596800247333334
$ws.Cells.Item(11,18).Value = 4.1371202226
$ws.Cells.Item(11,19).Value = 0.02157798815539813
$ws.Cells.Item(11,20).Value = 0.02157798815539813

# Row 12: Inflammatory-Mac -> Inflammatory-Mac
$ws.Cells.Item(12,1).Value = "Inflammatory-Mac"
$ws.Cells.Item(12,2).Value = "Nlgn2"
$ws.Cells.Item(12,3).Value = "Nrxn2"
$ws.Cells.Item(12,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(12,5).Value = 3
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = 0.9374000000000001
$ws.Cells.Item(12,8).Value = 2.8122
$ws.Cells.Item(12,9).Value = 0.03498421716936109
$ws.Cells.Item(12,10).Value = 0.03498421716936108
$ws.Cells.Item(12,11).Value = 1
$ws.Cells.Item(12,12).Value = 0.3333333333333333
$ws.Cells.Item(12,13).Value = 0.04417666666666667
$ws.Cells.Item(12,14).Value = 0.13253
$ws.Cells.Item(12,15).Value = 0.0555649656874912
$ws.Cells.Item(12,16).Value = 0.0555649656874912
$ws.Cells.Item(12,17).Value = 0.04141120733333334
$ws.Cells.Item(12,18).Value = 0.3727008660000001
$ws.Cells.Item(12,19).Value = 0.001943896826619289
$ws.Cells.Item(12,20).Value = 0.001943896826619289

# Row 13: Inflammatory-Mac -> MuSCs
$ws.Cells.Item(13,1).Value = "Inflammatory-Mac"
$ws.Cells.Item(13,2).Value = "Nlgn2"
$ws.Cells.Item(13,3).Value = "Nrxn2"
$ws.Cells.Item(13,4).Value = "MuSCs"
$ws.Cells.Item(13,5).Value = 3
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 0.9374000000000001
$ws.Cells.Item(13,8).Value = 2.8122
$ws.Cells.Item(13,9).Value = 0.03498421716936109
$ws.Cells.Item(13,10).Value = 0.03498421716936108
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 0.08703
$ws.Cells.Item(13,14).Value = 0.26109
$ws.Cells.Item(13,15).Value = 0.1094654560578516
$ws.Cells.Item(13,16).Value = 0.1094654560578516
$ws.Cells.Item(13,17).Value = 0.081581922
$ws.Cells.Item(13,18).Value = 0.734237298
$ws.Cells.Item(13,19).Value = 0.003829563287271034
$ws.Cells.Item(13,20).Value = 0.003829563287271034

# Row 14: MuSCs -> ECs
$ws.Cells.Item(14,1).Value = "MuSCs"
$ws.Cells.Item(14,2).Value = "Nlgn2"
$ws.Cells.Item(14,3).Value = "Nrxn2"
$ws.Cells.Item(14,4).Value = "ECs"
$ws.Cells.Item(14,5).Value = 3
$ws.Cells.Item(14,6).Value = 1
$ws.Cells.Item(14,7).Value = 3.238087
$ws.Cells.Item(14,8).Value = 9.714261
$ws.Cells.Item(14,9).Value = 0.1208469584182685
$ws.Cells.Item(14,10).Value = 0.1208469584182685
$ws.Cells.Item(14,11).Value = 3
$ws.Cells.Item(14,12).Value = 1
$ws.Cells.Item(14,13).Value = 0.173461
$ws.Cells.Item(14,14).Value = 0.520383
$ws.Cells.Item(14,15).Value = 0.2181774959583017
$ws.Cells.Item(14,16).Value = 0.2181774959583018
$ws.Cells.Item(14,17).Value = 0.5616818091070001
$ws.Cells.Item(14,18).Value = 5.055136281963001
$ws.Cells.Item(14,19).Value = 0.02636608678187484
$ws.Cells.Item(14,20).Value = 0.02636608678187485

# Row 15: MuSCs -> FAPs
$ws.Cells.Item(15,1).Value = "MuSCs"
$ws.Cells.Item(15,2).Value = "Nlgn2"
$ws.Cells.Item(15,3).Value = "Nrxn2"
$ws.Cells.Item(15,4).Value = "FAPs"
$ws.Cells.Item(15,5).Value = 3
$ws.Cells.Item(15,6).Value = 1
$ws.Cells.Item(15,7).Value = 3.238087
$ws.Cells.Item(15,8).Value = 9.714261
$ws.Cells.Item(15,9).Value = 0.1208469584182685
$ws.Cells.Item(15,10).Value = 0.1208469584182685
$ws.Cells.Item(15,11).Value = 3
$ws.Cells.Item(15,12).Value = 1
$ws.Cells.Item(15,13).Value = 0.4903776666666667
$ws.Cells.Item(15,14).Value = 1.471133
$ws.Cells.Item(15,15).Value = 0.6167920822963554
$ws.Cells.Item(15,16).Value = 0.6167920822963555
$ws.Cells.Item(15,17).Value = 1.587885547523667
$ws.Cells.Item(15,18).Value = 14.290969927713
$ws.Cells.Item(15,19).Value = 0.07453744712198491
$ws.Cells.Item(15,20).Value = 0.07453744712198493

# Row 16: MuSCs -> Inflammatory-Mac
$ws.Cells.Item(16,1).Value = "MuSCs"
$ws.Cells.Item(16,2).Value = "Nlgn2"
$ws.Cells.Item(16,3).Value = "Nrxn2"
$ws.Cells.Item(16,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(16,5).Value = 3
$ws.Cells.Item(16,6).Value = 1
$ws.Cells.Item(16,7).Value = 3.238087
$ws.Cells.Item(16,8).Value = 9.714261
$ws.Cells.Item(16,9).Value = 0.1208469584182685
$ws.Cells.Item(16,10).Value = 0.1208469584182685
$ws.Cells.Item(16,11).Value = 1
$ws.Cells.Item(16,12).Value = 0.3333333333333333
$ws.Cells.Item(16,13).Value = 0.04417666666666667
$ws.Cells.Item(16,14).Value = 0.13253
$ws.Cells.Item(16,15).Value = 0.0555649656874912
$ws.Cells.Item(16,16).Value = 0.0555649656874912
$ws.Cells.Item(16,17).Value = 0.1430478900366667
$ws.Cells.Item(16,18).Value = 1.28743101033
$ws.Cells.Item(16,19).Value = 0.006714857097948766
$ws.Cells.Item(16,20).Value = 0.006714857097948767

# Row 17: MuSCs -> MuSCs
$ws.Cells.Item(17,1).Value = "MuSCs"
$ws.Cells.Item(17,2).Value = "Nlgn2"
$ws.Cells.Item(17,3).Value = "Nrxn2"
$ws.Cells.Item(17,4).Value = "MuSCs"
$ws.Cells.Item(17,5).Value = 3
$ws.Cells.Item(17,6).Value = 1
$ws.Cells.Item(17,7).Value = 3.238087
$ws.Cells.Item(17,8).Value = 9.714261
$ws.Cells.Item(17,9).Value = 0.1208469584182685
$ws.Cells.Item(17,10).Value = 0.1208469584182685
$ws.Cells.Item(17,11).Value = 3
$ws.Cells.Item(17,12).Value = 1
$ws.Cells.Item(17,13).Value = 0.08703
$ws.Cells.Item(17,14).Value = 0.26109
$ws.Cells.Item(17,15).Value = 0.1094654560578516
$ws.Cells.Item(17,16).Value = 0.1094654560578516
$ws.Cells.Item(17,17).Value = 0.28181071161
$ws.Cells.Item(17,18).Value = 2.53629640449
$ws.Cells.Item(17,19).Value = 0.01322856741645999
$ws.Cells.Item(17,20).Value = 0.01322856741646

# Row 18: Resolving-Mac -> ECs
$ws.Cells.Item(18,1).Value = "Resolving-Mac"
$ws.Cells.Item(18,2).Value = "Nlgn2"
$ws.Cells.Item(18,3).Value = "Nrxn2"
$ws.Cells.Item(18,4).Value = "ECs"
$ws.Cells.Item(18,5).Value = 3
$ws.Cells.Item(18,6).Value = 1
$ws.Cells.Item(18,7).Value = 0.313797
$ws.Cells.Item(18,8).Value = 0.941391
$ws.Cells.Item(18,9).Value = 0.01171105440056966
$ws.Cells.Item(18,10).Value = 0.01171105440056966
$ws.Cells.Item(18,11).Value = 3
$ws.Cells.Item(18,12).Value = 1
$ws.Cells.Item(18,13).Value = 0.173461
$ws.Cells.Item(18,14).Value = 0.520383
$ws.Cells.Item(18,15).Value = 0.2181774959583017
$ws.Cells.Item(18,16).Value = 0.2181774959583018
$ws.Cells.Item(18,17).Value = 0.054431541417
$ws.Cells.Item(18,18).Value = 0.489883872753
$ws.Cells.Item(18,19).Value = 0.002555088524147739
$ws.Cells.Item(18,20).Value = 0.002555088524147739

# Row 19: Resolving-Mac -> FAPs
$ws.Cells.Item(19,1).Value = "Resolving-Mac"
$ws.Cells.Item(19,2).Value = "Nlgn2"
$ws.Cells.Item(19,3).Value = "Nrxn2"
$ws.Cells.Item(19,4).Value = "FAPs"
$ws.Cells.Item(19,5).Value = 3
$ws.Cells.Item(19,6).Value = 1
$ws.Cells.Item(19,7).Value = 0.313797
$ws.Cells.Item(19,8).Value = 0.941391
$ws.Cells.Item(19,9).Value = 0.01171105440056966
$ws.Cells.Item(19,10).Value = 0.01171105440056966
$ws.Cells.Item(19,11).Value = 3
$ws.Cells.Item(19,12).Value = 1
$ws.Cells.Item(19,13).Value = 0.4903776666666667
$ws.Cells.Item(19,14).Value = 1.471133
$ws.Cells.Item(19,15).Value = 0.6167920822963554
$ws.Cells.Item(19,16).Value = 0.6167920822963555
$ws.Cells.Item(19,17).Value = 0.153879040667
$ws.Cells.Item(19,18).Value = 1.384911366003
$ws.Cells.Item(19,19).Value = 0.007223285629613257
$ws.Cells.Item(19,20).Value = 0.007223285629613258

# Row 20: Resolving-Mac -> Inflammatory-Mac
$ws.Cells.Item(20,1).Value = "Resolving-Mac"
$ws.Cells.Item(20,2).Value = "Nlgn2"
$ws.Cells.Item(20,3).Value = "Nrxn2"
$ws.Cells.Item(20,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(20,5).Value = 3
$ws.Cells.Item(20,6).Value = 1
$ws.Cells.Item(20,7).Value = 0.313797
$ws.Cells.Item(20,8).Value = 0.941391
$ws.Cells.Item(20,9).Value = 0.01171105440056966
$ws.Cells.Item(20,10).Value = 0.01171105440056966
$ws.Cells.Item(20,11).Value = 1
$ws.Cells.Item(20,12).Value = 0.3333333333333333
$ws.Cells.Item(20,13).Value = 0.04417666666666667
$ws.Cells.Item(20,14).Value = 0.13253
$ws.Cells.Item(20,15).Value = 0.0555649656874912
$ws.Cells.Item(20,16).Value = 0.0555649656874912
$ws.Cells.Item(20,17).Value = 0.01386250547
$ws.Cells.Item(20,18).Value = 0.12476254923
$ws.Cells.Item(20,19).Value = 0.0006507243359319959
$ws.Cells.Item(20,20).Value = 0.000650724335931996

# Row 21: Resolving-Mac -> MuSCs
$ws.Cells.Item(21,1).Value = "Resolving-Mac"
$ws.Cells.Item(21,2).Value = "Nlgn2"
$ws.Cells.Item(21,3).Value = "Nrxn2"
$ws.Cells.Item(21,4).Value = "MuSCs"
$ws.Cells.Item(21,5).Value = 3
$ws.Cells.Item(21,6).Value = 1
$ws.Cells.Item(21,7).Value = 0.313797
$ws.Cells.Item(21,8).Value = 0.941391
$ws.Cells.Item(21,9).Value = 0.01171105440056966
$ws.Cells.Item(21,10).Value = 0.01171105440056966
$ws.Cells.Item(21,11).Value = 3
$ws.Cells.Item(21,12).Value = 1
$ws.Cells.Item(21,13).Value = 0.08703
$ws.Cells.Item(21,14).Value = 0.26109
$ws.Cells.Item(21,15).Value = 0.1094654560578516
$ws.Cells.Item(21,16).Value = 0.1094654560578516
$ws.Cells.Item(21,17).Value = 0.02730975291
$ws.Cells.Item(21,18).Value = 0.24578777619
$ws.Cells.Item(21,19).Value = 0.001281955910876668
$ws.Cells.Item(21,20).Value = 0.001281955910876668

Write-Output "Updated rows 2-21 (A:T) with refreshed TPM values"
